# Edit: Sat, Jul 11, 2020  6:04:51 PM
#
# 1) Slide 16: the table's style is changed to a different built-in
#    PowerPoint table style (GUID {568DCFCE-A950-44F3-9CF9-98B9CCB46A75}).
# 2) The deck's design/colour theme is changed from "Integral" to the
#    stock "Office Theme" palette (the two embedded theme parts end up
#    swapped: the palette that used to live in theme1.xml now lives in
#    theme2.xml, which is the one actually used by the slide master).

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 16 -------------------------------------------
$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{568DCFCE-A950-44F3-9CF9-98B9CCB46A75}")
    }
}

# --- 2. Swap the presentation's colour theme (Integral -> Office Theme) ---
# Table styles aside, Office Theme's twelve scheme colours replace the
# Integral ones that the slide master's theme currently carries.
$officeThemeRGB = @{
    1  = 0          # dk1      000000
    2  = 16777215   # lt1      FFFFFF
    3  = 6968388     # dk2      44546A
    4  = 15132391    # lt2      E7E6E6
    5  = 13998939    # accent1  5B9BD5
    6  = 3243501      # accent2  ED7D31
    7  = 10855845     # accent3  A5A5A5
    8  = 49407        # accent4  FFC000
    9  = 12874308     # accent5  4472C4
    10 = 4697456      # accent6  70AD47
    11 = 12673797     # hlink    0563C1
    12 = 7491477      # folHlink 954F72
}

$themeColors = $slide.ThemeColorScheme
for ($idx = 1; $idx -le 12; $idx++) {
    $themeColors.Colors($idx).RGB = $officeThemeRGB[$idx]
}
